$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7365
$ws.Range("B3").Value = 3035
$ws.Range("B4").Value = 6021
$ws.Range("B6").Value = 1730
$ws.Range("B7").Value = 108
$ws.Range("B8").Value = 116377
$ws.Range("B9").Value = 7278
$ws.Range("B10").Value = 1948
$ws.Range("B11").Value = 12015
$ws.Range("B12").Value = 7624
$ws.Range("B13").Value = 671
$ws.Range("B14").Value = 1394
$ws.Range("B15").Value = 27953
$ws.Range("B16").Value = 210
$ws.Range("B17").Value = 4948
$ws.Range("B18").Value = 26568
$ws.Range("B19").Value = 558
$ws.Range("B22").Value = 19062
$ws.Range("B23").Value = 12219
$ws.Range("B24").Value = 2416
$ws.Range("B25").Value = 612659
$ws.Range("B26").Value = 96
$ws.Range("B27").Value = 27253
$ws.Range("B28").Value = 265
$ws.Range("B29").Value = 19025
$ws.Range("B31").Value = 350
$ws.Range("B32").Value = 2900
$ws.Range("B33").Value = 1770
$ws.Range("B34").Value = 29550
$ws.Range("B35").Value = 101
$ws.Range("B37").Value = 38117
$ws.Range("B39").Value = 128093
$ws.Range("B40").Value = 150
$ws.Range("B41").Value = 339
$ws.Range("B42").Value = 1104
$ws.Range("B43").Value = 7251
$ws.Range("B44").Value = 702
$ws.Range("B45").Value = 10303
$ws.Range("B46").Value = 8295
$ws.Range("B47").Value = 590
$ws.Range("B48").Value = 32079
$ws.Range("B49").Value = 2822
$ws.Range("B50").Value = 186
$ws.Range("B51").Value = 35
$ws.Range("B52").Value = 4186
$ws.Range("B53").Value = 33121
$ws.Range("B54").Value = 19933
$ws.Range("B55").Value = 3753
$ws.Range("B56").Value = 170
$ws.Range("B57").Value = 54
$ws.Range("B58").Value = 1735
$ws.Range("B59").Value = 1248
$ws.Range("B60").Value = 6682
$ws.Range("B61").Value = 695
$ws.Range("B62").Value = 1249
$ws.Range("B63").Value = 119424
$ws.Range("B64").Value = 273
$ws.Range("B65").Value = 342
$ws.Range("B66").Value = 11498
$ws.Range("B67").Value = 99130
$ws.Range("B68").Value = 1208
$ws.Range("B69").Value = 17313
$ws.Range("B70").Value = 200
$ws.Range("B71").Value = 15840
$ws.Range("B72").Value = 387
$ws.Range("B73").Value = 146
$ws.Range("B74").Value = 976
$ws.Range("B75").Value = 711
$ws.Range("B76").Value = 10387
$ws.Range("B77").Value = 32780
$ws.Range("B79").Value = 465911
$ws.Range("B80").Value = 143739
$ws.Range("B81").Value = 128956
$ws.Range("B82").Value = 23665
$ws.Range("B83").Value = 5609
$ws.Range("B84").Value = 8176
$ws.Range("B85").Value = 133177
$ws.Range("B86").Value = 2358
$ws.Range("B87").Value = 18343
$ws.Range("B88").Value = 11361
$ws.Range("B89").Value = 17625
$ws.Range("B90").Value = 5328
$ws.Range("B92").Value = 3298
$ws.Range("B93").Value = 2982
$ws.Range("B94").Value = 2464
$ws.Range("B95").Value = 2723
$ws.Range("B96").Value = 132
$ws.Range("B97").Value = 3950
$ws.Range("B98").Value = 8645
$ws.Range("B99").Value = 661
$ws.Range("B101").Value = 5380
$ws.Range("B103").Value = 6521
$ws.Range("B104").Value = 862
$ws.Range("B106").Value = 2304
$ws.Range("B107").Value = 30002
$ws.Range("B108").Value = 248
$ws.Range("B109").Value = 592
$ws.Range("B110").Value = 463
$ws.Range("B112").Value = 818
$ws.Range("B113").Value = 240
$ws.Range("B114").Value = 292372
$ws.Range("B115").Value = 8834
$ws.Range("B117").Value = 1948
$ws.Range("B118").Value = 2245
$ws.Range("B119").Value = 14761
$ws.Range("B120").Value = 1938
$ws.Range("B121").Value = 3569
$ws.Range("B122").Value = 11496
$ws.Range("B123").Value = 19411
$ws.Range("B124").Value = 40
$ws.Range("B125").Value = 211
$ws.Range("B126").Value = 242
$ws.Range("B127").Value = 2974
$ws.Range("B128").Value = 999
$ws.Range("B129").Value = 4113
$ws.Range("B130").Value = 28663
$ws.Range("B132").Value = 7353
$ws.Range("B133").Value = 486
$ws.Range("B134").Value = 16354
$ws.Range("B135").Value = 200866
$ws.Range("B136").Value = 47074
$ws.Range("B137").Value = 80822
$ws.Range("B138").Value = 18321
$ws.Range("B140").Value = 54959
$ws.Range("B141").Value = 259107
$ws.Range("B142").Value = 1340
$ws.Range("B143").Value = 28
$ws.Range("B144").Value = 278
$ws.Range("B145").Value = 73
$ws.Range("B147").Value = 93
$ws.Range("B149").Value = 8824
$ws.Range("B150").Value = 1883
$ws.Range("B151").Value = 11225
$ws.Range("B152").Value = 125
$ws.Range("B154").Value = 662
$ws.Range("B155").Value = 13861
$ws.Range("B156").Value = 5063
$ws.Range("B158").Value = 1324
$ws.Range("B159").Value = 89574
$ws.Range("B161").Value = 87810
$ws.Range("B162").Value = 14127
$ws.Range("B164").Value = 1149
$ws.Range("B165").Value = 15110
$ws.Range("B166").Value = 11382
$ws.Range("B167").Value = 2704
$ws.Range("B168").Value = 848
$ws.Range("B170").Value = 727
$ws.Range("B171").Value = 20387
$ws.Range("B175").Value = 1967
$ws.Range("B176").Value = 25344
$ws.Range("B177").Value = 75042
$ws.Range("B178").Value = 771118
$ws.Range("B179").Value = 3256
$ws.Range("B180").Value = 86260
$ws.Range("B181").Value = 2144
$ws.Range("B182").Value = 144369
$ws.Range("B183").Value = 6115
$ws.Range("B184").Value = 1375
$ws.Range("B186").Value = 5079
$ws.Range("B187").Value = 23761
$ws.Range("B188").Value = 4770
$ws.Range("B189").Value = 1938
$ws.Range("B190").Value = 3667
$ws.Range("B191").Value = 4699
